$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-coerced to a number by Excel
# (e.g. "1.00", "219.48") are forced back to Text format first so the stored
# cell keeps the exact literal string shown in the source diff.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D16", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D39", "D41", "D44", "D45", "D46", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "89.385.34"
$ws.Range("E2").Value = "  +10.24%  "
$ws.Range("D3").Value = "3.361.83"
$ws.Range("E3").Value = "  +6.96%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "219.48"
$ws.Range("E5").Value = "  +5.43%  "
$ws.Range("D6").Value = "647.25"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("D7").Value = "0.415"
$ws.Range("E7").Value = "  +47.64%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.613"
$ws.Range("E9").Value = "  +6.43%  "
$ws.Range("D10").Value = "3.360.80"
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("D11").Value = "0.618"
$ws.Range("E11").Value = "  +7.75%  "
$ws.Range("D12").Value = "0.0000291"
$ws.Range("E12").Value = "  +16.08%  "
$ws.Range("D13").Value = "36.42"
$ws.Range("E13").Value = "  +16.11%  "
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "3.977.51"
$ws.Range("E15").Value = "  +7.02%  "
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("D17").Value = "89.216.27"
$ws.Range("E17").Value = "  +10.38%  "
$ws.Range("D18").Value = "3.365.48"
$ws.Range("E18").Value = "  +7.38%  "
$ws.Range("D19").Value = "14.82"
$ws.Range("E19").Value = "  +7.15%  "
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "9.71"
$ws.Range("E21").Value = "  +8.65%  "
$ws.Range("D22").Value = "458.31"
$ws.Range("E22").Value = "  +6.71%  "
$ws.Range("D23").Value = "5.56"
$ws.Range("E23").Value = "  +9.67%  "
$ws.Range("D24").Value = "7.49"
$ws.Range("E24").Value = "  +4.65%  "
$ws.Range("D25").Value = "5.56"
$ws.Range("E25").Value = "  +7.20%  "
$ws.Range("D26").Value = "12.78"
$ws.Range("E26").Value = "  +18.39%  "
$ws.Range("D27").Value = "3.524.32"
$ws.Range("E27").Value = "  +6.85%  "
$ws.Range("D28").Value = "0.0000144"
$ws.Range("E28").Value = "  +19.10%  "
$ws.Range("D29").Value = "79.36"
$ws.Range("E29").Value = "  +4.96%  "
$ws.Range("D30").Value = "0.199"
$ws.Range("E30").Value = "  +45.96%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "9.42"
$ws.Range("E32").Value = "  +5.60%  "
$ws.Range("D33").Value = "595.06"
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E35").Value = "  +7.02%  "
$ws.Range("E36").Value = "  +6.98%  "
$ws.Range("D37").Value = "7.30"
$ws.Range("E37").Value = "  +19.98%  "
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").Value = "23.55"
$ws.Range("E39").Value = "  +4.59%  "
$ws.Range("E40").Value = "  +6.26%  "
$ws.Range("D41").Value = "2.16"
$ws.Range("E41").Value = "  +7.13%  "
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("E43").Value = "  +5.33%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").Value = "1.46"
$ws.Range("E45").Value = "  +11.12%  "
$ws.Range("D46").Value = "157.97"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "188.89"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "46.44"
$ws.Range("E49").Value = "  +4.08%  "
$ws.Range("E50").Value = "  +7.94%  "
$ws.Range("E51").Value = "  +7.14%  "
